$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Columns.Item(20).Insert()
$ws.Rows.Item(13).Insert()
$ws.Cells.Item(1, 20).Value = "media::image"
$ws.Cells.Item(13, 1).Value = "note"
$ws.Cells.Item(13, 2).Value = "logo"
$ws.Cells.Item(13, 20).Value = "logo.png"
